$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).Value = 'Record'
$ws.Cells.Item($row, 2).Value = 'RJ Record'
$ws.Cells.Item($row, 3).Value = 'Defesa Civil'
$ws.Cells.Item($row, 4).Value = '2025-04-03T18:15'
$ws.Cells.Item($row, 5).Value = 'Positivo'
$ws.Cells.Item($row, 6).Value = 'Vem chuva! Nova frente fria se aproxima do Estado e a previsão é de temporais. Repórter *ao vivo*. Em Campos, Centro de Monitoramento de Desastres da Defesa Civil está atento à situação. entrevista com o coordenador do Centro, Leandro Freitas, que orientou sobre o sistema de alertas. Enviar CEP para 40199 por SMS.  '
